$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.747155
$ws.Range("H2").Value = 23.241465
$ws.Range("I2").Value = 0.1917182689216984
$ws.Range("J2").Value = 0.2322720883603271
$ws.Range("M2").Value = 54.23134333333334
$ws.Range("N2").Value = 162.69403
$ws.Range("O2").Value = 0.9097185042023884
$ws.Range("P2").Value = 0.9200351849746305
$ws.Range("Q2").Value = 420.13862266155
$ws.Range("R2").Value = 3781.24760395395
$ws.Range("S2").Value = 0.1744096568317187
$ws.Range("T2").Value = 0.2136984937790373
$ws.Range("G3").Value = 7.747155
$ws.Range("H3").Value = 23.241465
$ws.Range("I3").Value = 0.1917182689216984
$ws.Range("J3").Value = 0.2322720883603271
$ws.Range("O3").Value = 0.04961247197704183
$ws.Range("P3").Value = 0.0501751032012552
$ws.Range("Q3").Value = 22.912709312805
$ws.Range("R3").Value = 206.214383815245
$ws.Range("S3").Value = 0.009511617244364728
$ws.Range("T3").Value = 0.01165427600425048
$ws.Range("G4").Value = 7.747155
$ws.Range("H4").Value = 23.241465
$ws.Range("I4").Value = 0.1917182689216984
$ws.Range("J4").Value = 0.2322720883603271
$ws.Range("M4").Value = 0.2447093333333333
$ws.Range("N4").Value = 0.734128
$ws.Range("O4").Value = 0.004104943654374356
$ws.Range("P4").Value = 0.004151495849448536
$ws.Range("Q4").Value = 1.89580113528
$ws.Range("R4").Value = 17.06221021752
$ws.Range("S4").Value = 0.0007869926914377619
$ws.Range("T4").Value = 0.0009642766107706414
$ws.Range("G5").Value = 7.747155
$ws.Range("H5").Value = 23.241465
$ws.Range("I5").Value = 0.1917182689216984
$ws.Range("J5").Value = 0.2322720883603271
$ws.Range("M5").Value = 2.005396
$ws.Range("N5").Value = 4.010792
$ws.Range("O5").Value = 0.03364006379558217
$ws.Range("P5").Value = 0.02268103973830366
$ws.Range("Q5").Value = 15.53611364838
$ws.Range("R5").Value = 93.21668189028001
$ws.Range("S5").Value = 0.00644941479730451
$ws.Range("T5").Value = 0.005268172466199358
$ws.Range("G6").Value = 7.747155
$ws.Range("H6").Value = 23.241465
$ws.Range("I6").Value = 0.1917182689216984
$ws.Range("J6").Value = 0.2322720883603271
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1743103333333333
$ws.Range("N6").Value = 0.522931
$ws.Range("O6").Value = 0.002924016370613348
$ws.Range("P6").Value = 0.002957176236362014
$ws.Range("Q6").Value = 1.350409170435
$ws.Range("R6").Value = 12.153682533915
$ws.Range("S6").Value = 0.0005605873568726984
$ws.Range("T6").Value = 0.0006868695000693372
$ws.Range("I7").Value = 0.2330460307376861
$ws.Range("J7").Value = 0.2823418370506734
$ws.Range("M7").Value = 54.23134333333334
$ws.Range("N7").Value = 162.69403
$ws.Range("O7").Value = 0.9097185042023884
$ws.Range("P7").Value = 0.9200351849746305
$ws.Range("Q7").Value = 510.70583372971
$ws.Range("R7").Value = 4596.35250356739
$ws.Range("S7").Value = 0.2120062864929917
$ws.Range("T7").Value = 0.2597644242769933
$ws.Range("I8").Value = 0.2330460307376861
$ws.Range("J8").Value = 0.2823418370506734
$ws.Range("O8").Value = 0.04961247197704183
$ws.Range("P8").Value = 0.0501751032012552
$ws.Range("S8").Value = 0.01156198966933428
$ws.Range("T8").Value = 0.01416653081204952
$ws.Range("I9").Value = 0.2330460307376861
$ws.Range("J9").Value = 0.2823418370506734
$ws.Range("M9").Value = 0.2447093333333333
$ws.Range("N9").Value = 0.734128
$ws.Range("O9").Value = 0.004104943654374356
$ws.Range("P9").Value = 0.004151495849448536
$ws.Range("Q9").Value = 2.304469637296
$ws.Range("R9").Value = 20.740226735664
$ws.Range("S9").Value = 0.0009566408250537958
$ws.Range("T9").Value = 0.001172140964641545
$ws.Range("I10").Value = 0.2330460307376861
$ws.Range("J10").Value = 0.2823418370506734
$ws.Range("M10").Value = 2.005396
$ws.Range("N10").Value = 4.010792
$ws.Range("O10").Value = 0.03364006379558217
$ws.Range("P10").Value = 0.02268103973830366
$ws.Range("Q10").Value = 18.885157054716
$ws.Range("R10").Value = 113.310942328296
$ws.Range("S10").Value = 0.007839683341322965
$ws.Range("T10").Value = 0.006403806425931981
$ws.Range("I11").Value = 0.2330460307376861
$ws.Range("J11").Value = 0.2823418370506734
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1743103333333333
$ws.Range("N11").Value = 0.522931
$ws.Range("O11").Value = 0.002924016370613348
$ws.Range("P11").Value = 0.002957176236362014
$ws.Range("Q11").Value = 1.641510216067
$ws.Range("R11").Value = 14.773591944603
$ws.Range("S11").Value = 0.0006814304089834558
$ws.Range("T11").Value = 0.0008349345710570473
$ws.Range("G12").Value = 1.149447666666667
$ws.Range("H12").Value = 3.448343
$ws.Range("I12").Value = 0.02844529596599251
$ws.Range("J12").Value = 0.03446227808757819
$ws.Range("M12").Value = 54.23134333333334
$ws.Range("N12").Value = 162.69403
$ws.Range("O12").Value = 0.9097185042023884
$ws.Range("P12").Value = 0.9200351849746305
$ws.Range("Q12").Value = 62.33609105469889
$ws.Range("R12").Value = 561.02481949229
$ws.Range("S12").Value = 0.02587721209777694
$ws.Range("T12").Value = 0.03170650839495216
$ws.Range("G13").Value = 1.149447666666667
$ws.Range("H13").Value = 3.448343
$ws.Range("I13").Value = 0.02844529596599251
$ws.Range("J13").Value = 0.03446227808757819
$ws.Range("O13").Value = 0.04961247197704183
$ws.Range("P13").Value = 0.0501751032012552
$ws.Range("Q13").Value = 3.399565421966555
$ws.Range("R13").Value = 30.596088797699
$ws.Range("S13").Value = 0.001411241448991464
$ws.Range("T13").Value = 0.001729148359594591
$ws.Range("G14").Value = 1.149447666666667
$ws.Range("H14").Value = 3.448343
$ws.Range("I14").Value = 0.02844529596599251
$ws.Range("J14").Value = 0.03446227808757819
$ws.Range("M14").Value = 0.2447093333333333
$ws.Range("N14").Value = 0.734128
$ws.Range("O14").Value = 0.004104943654374356
$ws.Range("P14").Value = 0.004151495849448536
$ws.Range("Q14").Value = 0.2812805722115556
$ws.Range("R14").Value = 2.531525149904
$ws.Range("S14").Value = 0.0001167663371724014
$ws.Range("T14").Value = 0.0001430700044431221
$ws.Range("G15").Value = 1.149447666666667
$ws.Range("H15").Value = 3.448343
$ws.Range("I15").Value = 0.02844529596599251
$ws.Range("J15").Value = 0.03446227808757819
$ws.Range("M15").Value = 2.005396
$ws.Range("N15").Value = 4.010792
$ws.Range("O15").Value = 0.03364006379558217
$ws.Range("P15").Value = 0.02268103973830366
$ws.Range("Q15").Value = 2.305097752942667
$ws.Range("R15").Value = 13.830586517656
$ws.Range("S15").Value = 0.000956901570980204
$ws.Range("T15").Value = 0.0007816402987768325
$ws.Range("G16").Value = 1.149447666666667
$ws.Range("H16").Value = 3.448343
$ws.Range("I16").Value = 0.02844529596599251
$ws.Range("J16").Value = 0.03446227808757819
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1743103333333333
$ws.Range("N16").Value = 0.522931
$ws.Range("O16").Value = 0.002924016370613348
$ws.Range("P16").Value = 0.002957176236362014
$ws.Range("Q16").Value = 0.2003606059258889
$ws.Range("R16").Value = 1.803245453333
$ws.Range("S16").Value = 0.00008317451107150392
$ws.Range("T16").Value = 0.0001019110298114856
$ws.Range("G17").Value = 21.1658025
$ws.Range("H17").Value = 42.331605
$ws.Range("I17").Value = 0.523788541166216
$ws.Range("J17").Value = 0.4230563906790929
$ws.Range("M17").Value = 54.23134333333334
$ws.Range("N17").Value = 162.69403
$ws.Range("O17").Value = 0.9097185042023884
$ws.Range("P17").Value = 0.9200351849746305
$ws.Range("Q17").Value = 1147.849902303025
$ws.Range("R17").Value = 6887.09941381815
$ws.Range("S17").Value = 0.4765001281880811
$ws.Range("T17").Value = 0.3892267646531388
$ws.Range("G18").Value = 21.1658025
$ws.Range("H18").Value = 42.331605
$ws.Range("I18").Value = 0.523788541166216
$ws.Range("J18").Value = 0.4230563906790929
$ws.Range("O18").Value = 0.04961247197704183
$ws.Range("P18").Value = 0.0501751032012552
$ws.Range("Q18").Value = 62.59922256037751
$ws.Range("R18").Value = 375.595335362265
$ws.Range("S18").Value = 0.02598644432050451
$ws.Range("T18").Value = 0.02122689806227402
$ws.Range("G19").Value = 21.1658025
$ws.Range("H19").Value = 42.331605
$ws.Range("I19").Value = 0.523788541166216
$ws.Range("J19").Value = 0.4230563906790929
$ws.Range("M19").Value = 0.2447093333333333
$ws.Range("N19").Value = 0.734128
$ws.Range("O19").Value = 0.004104943654374356
$ws.Range("P19").Value = 0.004151495849448536
$ws.Range("Q19").Value = 5.17946941924
$ws.Range("R19").Value = 31.07681651544
$ws.Range("S19").Value = 0.002150122448294259
$ws.Range("T19").Value = 0.001756316849986932
$ws.Range("G20").Value = 21.1658025
$ws.Range("H20").Value = 42.331605
$ws.Range("I20").Value = 0.523788541166216
$ws.Range("J20").Value = 0.4230563906790929
$ws.Range("M20").Value = 2.005396
$ws.Range("N20").Value = 4.010792
$ws.Range("O20").Value = 0.03364006379558217
$ws.Range("P20").Value = 0.02268103973830366
$ws.Range("Q20").Value = 42.44581567029001
$ws.Range("R20").Value = 169.78326268116
$ws.Range("S20").Value = 0.01762027994022642
$ws.Range("T20").Value = 0.009595358808535827
$ws.Range("G21").Value = 21.1658025
$ws.Range("H21").Value = 42.331605
$ws.Range("I21").Value = 0.523788541166216
$ws.Range("J21").Value = 0.4230563906790929
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.1743103333333333
$ws.Range("N21").Value = 0.522931
$ws.Range("O21").Value = 0.002924016370613348
$ws.Range("P21").Value = 0.002957176236362014
$ws.Range("Q21").Value = 3.6894180890425
$ws.Range("R21").Value = 22.136508534255
$ws.Range("S21").Value = 0.001531566269109699
$ws.Range("T21").Value = 0.001251052305157298
$ws.Range("G22").Value = 0.9294836666666667
$ws.Range("H22").Value = 2.788451
$ws.Range("I22").Value = 0.02300186320840699
$ws.Range("J22").Value = 0.02786740582232843
$ws.Range("M22").Value = 54.23134333333334
$ws.Range("N22").Value = 162.69403
$ws.Range("O22").Value = 0.9097185042023884
$ws.Range("P22").Value = 0.9200351849746305
$ws.Range("Q22").Value = 50.40714784972556
$ws.Range("R22").Value = 453.6643306475301
$ws.Range("S22").Value = 0.02092522059181996
$ws.Range("T22").Value = 0.02563899387050903
$ws.Range("G23").Value = 0.9294836666666667
$ws.Range("H23").Value = 2.788451
$ws.Range("I23").Value = 0.02300186320840699
$ws.Range("J23").Value = 0.02786740582232843
$ws.Range("O23").Value = 0.04961247197704183
$ws.Range("P23").Value = 0.0501751032012552
$ws.Range("Q23").Value = 2.749007740949223
$ws.Range("R23").Value = 24.741069668543
$ws.Range("S23").Value = 0.001141179293846841
$ws.Range("T23").Value = 0.001398249963086589
$ws.Range("G24").Value = 0.9294836666666667
$ws.Range("H24").Value = 2.788451
$ws.Range("I24").Value = 0.02300186320840699
$ws.Range("J24").Value = 0.02786740582232843
$ws.Range("M24").Value = 0.2447093333333333
$ws.Range("N24").Value = 0.734128
$ws.Range("O24").Value = 0.004104943654374356
$ws.Range("P24").Value = 0.004151495849448536
$ws.Range("Q24").Value = 0.2274533284142222
$ws.Range("R24").Value = 2.047079955728
$ws.Range("S24").Value = 0.00009442135241613722
$ws.Range("T24").Value = 0.0001156914196062944
$ws.Range("G25").Value = 0.9294836666666667
$ws.Range("H25").Value = 2.788451
$ws.Range("I25").Value = 0.02300186320840699
$ws.Range("J25").Value = 0.02786740582232843
$ws.Range("M25").Value = 2.005396
$ws.Range("N25").Value = 4.010792
$ws.Range("O25").Value = 0.03364006379558217
$ws.Range("P25").Value = 0.02268103973830366
$ws.Range("Q25").Value = 1.863982827198667
$ws.Range("R25").Value = 11.183896963192
$ws.Range("S25").Value = 0.0007737841457480653
$ws.Range("T25").Value = 0.000632061738859666
$ws.Range("G26").Value = 0.9294836666666667
$ws.Range("H26").Value = 2.788451
$ws.Range("I26").Value = 0.02300186320840699
$ws.Range("J26").Value = 0.02786740582232843
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.1743103333333333
$ws.Range("N26").Value = 0.522931
$ws.Range("O26").Value = 0.002924016370613348
$ws.Range("P26").Value = 0.002957176236362014
$ws.Range("Q26").Value = 0.1620186077645556
$ws.Range("R26").Value = 1.458167469881
$ws.Range("S26").Value = 0.00006725782457599091
$ws.Range("T26").Value = 0.00008240883026684605

Write-Host "Updated 288 cells"
